$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Precondition text (row 8 and row 18 share the same text originally)
$ws.Range("B8").Value = "Lider de Pessoas esta autenticado no sistema e tem permissao para gerenciar Metas de Desempenho"
$ws.Range("B18").Value = "Lider de Pessoas esta autenticado no sistema e tem permissao para gerenciar Metas de Desempenho"

# Step 2 text in first test case (row 11)
$ws.Range("B11").Value = "Lider de Pessoas com uma avaliacao selecionada, clica na opcao 'Editar' para modificar a Avaliacao de Desempenho"
$ws.Range("D11").Value = "SYSTEM apresenta o formulario com o campo 'Metas' contendo cada Competencia do perfil avaliado"
